# Start Arthur Player Controller
#
# The "Arthur" tasks that lived in the Wednesday column (J) get carried over
# into the Thursday column (K); the Wednesday column is left empty. A new
# Thursday task is added ("Pick between physics and direct transform") and
# what used to be the "Make Plan For Thursday" note becomes "Make Plan For
# Friday", now also living in column K. Finally, the active selection moves
# from J6 to H3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: J2 ("Pseudocode out Arthur Jump") -> K2 ------------------------
$ws.Range("J2").Cut($ws.Range("K2"))
$ws.Range("J2").Clear()

# --- Row 3: J3 ("Implement Arthur Jump") -> K3 ------------------------------
# Row 3 also shrinks from a 60pt row to 45pt now that it only needs to fit
# the shorter remaining entries.
$ws.Range("J3").Cut($ws.Range("K3"))
$ws.Range("J3").Clear()
$ws.Range("3:3").RowHeight = 45

# --- Row 4: J4 ("Implement Arthur Sprite") -> K4 ----------------------------
$ws.Range("J4").Cut($ws.Range("K4"))
$ws.Range("J4").Clear()

# --- Row 5: new task added in K5, matching the centered/wrapped style ------
# used by the other header-ish cells in this column (same look as K2).
$ws.Range("K5").Value2 = "Pick between physics and direct transform"
$ws.Range("K5").HorizontalAlignment = -4108   # xlCenter
$ws.Range("K5").WrapText = $true

# --- Row 6: the old "Make Plan For Thursday" note becomes the new ----------
# "Make Plan For Friday" note, and moves from J6 to K6.
$ws.Range("J6").Value2 = "Make Plan For Friday"
$ws.Range("J6").Cut($ws.Range("K6"))
$ws.Range("J6").Clear()

# --- Update the active selection from J6 to H3 ------------------------------
[void]$ws.Range("H3").Select()
